# Saldo_guide.xlsx update: advance the reference date by one day
# (2024-09-18 -> 2024-09-19, Excel serial 45553 -> 45554) and refresh the
# balances that shifted for that new reference date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new extraction timestamp.
$ws.Name = "IClientBalance-20240919-092647-"

# Column G ("Dt. Referencia") holds the same date serial (45553) for every
# data row (2..274). Bump every one of them to 45554.
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45554
}

# A handful of rows also got revised "Saldo Previsto" / "Vl. Total" amounts
# (columns E and H) alongside the date change.
$ws.Range("E57").Value = 38072.699999999997
$ws.Range("H57").Value = 38072.699999999997

$ws.Range("E103").Value = 101070.12
$ws.Range("H103").Value = 101070.12

$ws.Range("E104").Value = -386.13
$ws.Range("H104").Value = -386.13

$ws.Range("E148").Value = 27791.79
$ws.Range("H148").Value = 27791.79

$ws.Range("E165").Value = 12582.81
$ws.Range("H165").Value = 12582.81
